$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 40 (shifts existing rows 40-124 down to 41-125)
$ws.Rows.Item(40).Insert()

# Populate the new row with the IC9 part data
$ws.Range("A40").Value = "IC9"
$ws.Range("C40").Value = "742G17DCK"
$ws.Range("D40").Value = "DCK_R-PDSO-G6"
$ws.Range("E40").Value = "ME"
$ws.Range("F40").Value = "595-SN74AUP2G17DCKR"
$ws.Range("G40").Value = "SN74AUP2G17DCKR"
$ws.Range("H40").Value = 1
$ws.Range("I40").Value = 0.47
$ws.Range("J40").Formula = "=H40*I40"

# Update the view state to match the saved selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("J40").Select()
